$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New continuation of the monthly date series for rows 338-560 (column A),
# replacing the erroneous jumped/duplicated dates that used the wrong date style.
$dates = @(36892,36923,36951,36982,37012,37043,37073,37104,37135,37165,37196,37226,37257,37288,37316,37347,37377,37408,37438,37469,37500,37530,37561,37591,37622,37653,37681,37712,37742,37773,37803,37834,37865,37895,37926,37956,37987,38018,38047,38078,38108,38139,38169,38200,38231,38261,38292,38322,38353,38384,38412,38443,38473,38504,38534,38565,38596,38626,38657,38687,38718,38749,38777,38808,38838,38869,38899,38930,38961,38991,39022,39052,39083,39114,39142,39173,39203,39234,39264,39295,39326,39356,39387,39417,39448,39479,39508,39539,39569,39600,39630,39661,39692,39722,39753,39783,39814,39845,39873,39904,39934,39965,39995,40026,40057,40087,40118,40148,40179,40210,40238,40269,40299,40330,40360,40391,40422,40452,40483,40513,40544,40575,40603,40634,40664,40695,40725,40756,40787,40817,40848,40878,40909,40940,40969,41000,41030,41061,41091,41122,41153,41183,41214,41244,41275,41306,41334,41365,41395,41426,41456,41487,41518,41548,41579,41609,41640,41671,41699,41730,41760,41791,41821,41852,41883,41913,41944,41974,42005,42036,42064,42095,42125,42156,42186,42217,42248,42278,42309,42339,42370,42401,42430,42461,42491,42522,42552,42583,42614,42644,42675,42705,42736,42767,42795,42826,42856,42887,42917,42948,42979,43009,43040,43070,43101,43132,43160,43191,43221,43252,43282,43313,43344,43374,43405,43435,43466,43497,43525,43556,43586,43617,43647)

$arr = New-Object "object[,]" 223,1
for ($i = 0; $i -lt 223; $i++) {
    $arr[$i,0] = $dates[$i]
}

$rng = $ws.Range("A338:A560")
$rng.Value = $arr
$rng.NumberFormat = "mmm-yy"

# Restore the sheet view / selection state recorded after the edit:
# the user had scrolled down and selected the corrected date range.
$ws.Range("A338:A560").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 616
$win.ScrollColumn = 1

